$d = $word.ActiveDocument

# Locate the paragraph that ends with the "Entidade(s) afetada(s)..." sentence.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Entidade(s) afetada(s) pelo rompimento do cabo de fibras óptica:*") {
        $target = $p
    }
}

if ($target -ne $null) {
    $anchor = $target.Range
    $anchor.Collapse(0)

    # Blank spacer paragraph (no pStyle, same spacing as surrounding body text).
    $spacer = $anchor.Paragraphs.Add($anchor)
    $spacer.Range.Text = ""
    $spacer.Format.LineSpacingRule = 5
    $spacer.Format.LineSpacing = $d.Application.LinesToPoints(1)
    $spacer.SpaceBefore = 0
    $spacer.SpaceAfter = 0
    $anchor = $spacer.Range
    $anchor.Collapse(0)

    $schools = @("Escola 01", "Escola 02", "Escola 03")
    foreach ($school in $schools) {
        $para = $anchor.Paragraphs.Add($anchor)
        $pRange = $para.Range
        $para.Style = "ListBullet"
        $para.Format.LineSpacingRule = 5
        $para.Format.LineSpacing = $d.Application.LinesToPoints(1)
        $para.SpaceBefore = 0
        $para.SpaceAfter = 0
        $para.LeftIndent = $d.Application.InchesToPoints(0.5)
        $para.Alignment = 0

        $pRange.Text = $school
        $pRange.Font.Name = "Arial"
        $pRange.Font.Bold = $true
        $pRange.Font.Italic = $false
        $pRange.Font.Size = 12

        $anchor = $para.Range
        $anchor.Collapse(0)
    }
}
